$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 609.8570999999999
$ws.Range("I2").Value = 660
$ws.Range("K2").Value = 660
$ws.Range("M2").Value = -547
$ws.Range("H41").Value = 298.57144
$ws.Range("I41").Value = 340
$ws.Range("J41").Value = 291.66666
$ws.Range("K41").Value = 340
$ws.Range("L41").Value = 291.66666
$ws.Range("M41").Value = 100
$ws.Range("N41").Value = -1171.66666
$ws.Range("H76").Value = 7999.8335
$ws.Range("J76").Value = 7999.8335
$ws.Range("L76").Value = 7999.8335
$ws.Range("N76").Value = -8629.833500000001
$ws.Range("H79").Value = 7999.8335
$ws.Range("J79").Value = 7999.8335
$ws.Range("L79").Value = 7999.8335
$ws.Range("N79").Value = -10183.8335
$ws.Range("H87").Value = 89999
$ws.Range("J87").Value = 89999
$ws.Range("L87").Value = 89999
$ws.Range("N87").Value = -92495
$ws.Range("H90").Value = 89999
$ws.Range("J90").Value = 89999
$ws.Range("L90").Value = 269997
$ws.Range("N90").Value = -282477
$ws.Range("H100").Value = 586.1905
$ws.Range("I100").Value = 518.82355
$ws.Range("K100").Value = 518.82355
$ws.Range("M100").Value = 22.17645000000005
$ws.Range("H112").Value = 2619.158
$ws.Range("J112").Value = 2809.3125
$ws.Range("L112").Value = 8427.9375
$ws.Range("N112").Value = -10643.9375
$ws.Range("H137").Value = 2794.6128
$ws.Range("I137").Value = 1530.2667
$ws.Range("J137").Value = 3979.9375
$ws.Range("K137").Value = 4590.800099999999
$ws.Range("L137").Value = 11939.8125
$ws.Range("M137").Value = -2040.800099999999
$ws.Range("N137").Value = -17039.8125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1491.8572
$ws.Range("I61").Value = 1491.8572
$ws.Range("K61").Value = 1491.8572
$ws.Range("M61").Value = -1279.8572
$ws.Range("H132").Value = 2312.4707
$ws.Range("I132").Value = 2312.4707
$ws.Range("K132").Value = 6937.4121
$ws.Range("M132").Value = -4407.4121
$ws.Range("H135").Value = 199999
$ws.Range("J135").Value = 199999
$ws.Range("L135").Value = 199999
$ws.Range("N135").Value = -210139
$ws.Range("H136").Value = 1491.8572
$ws.Range("I136").Value = 1491.8572
$ws.Range("K136").Value = 4475.571599999999
$ws.Range("M136").Value = -1925.571599999999
$ws.Range("H139").Value = 65499.25
$ws.Range("I139").Value = 54999
$ws.Range("J139").Value = 75999.5
$ws.Range("K139").Value = 54999
$ws.Range("L139").Value = 75999.5
$ws.Range("M139").Value = -49859
$ws.Range("N139").Value = -86279.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1599.5
$ws.Range("I86").Value = 1628.1428
$ws.Range("J86").Value = 1399
$ws.Range("K86").Value = 1628.1428
$ws.Range("L86").Value = 1399
$ws.Range("M86").Value = -505.1428000000001
$ws.Range("N86").Value = -3645
$ws.Range("H89").Value = 1599.5
$ws.Range("I89").Value = 1628.1428
$ws.Range("J89").Value = 1399
$ws.Range("K89").Value = 8140.714
$ws.Range("L89").Value = 6995
$ws.Range("M89").Value = -2524.714
$ws.Range("N89").Value = -18227
$ws.Range("H134").Value = 1812.2
$ws.Range("I134").Value = 765.5
$ws.Range("K134").Value = 2296.5
$ws.Range("M134").Value = 238.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 25001400
$ws.Range("I2").Value = 28572814
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 28572814
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -28572701
$ws.Range("N2").Value = -1726
$ws.Range("H58").Value = 7954.2
$ws.Range("I58").Value = 6000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5797
$ws.Range("H132").Value = 3977.818
$ws.Range("I132").Value = 3203.9333
$ws.Range("K132").Value = 9611.7999
$ws.Range("M132").Value = -7081.7999
$ws.Range("H136").Value = 7954.2
$ws.Range("I136").Value = 6000
$ws.Range("K136").Value = 18000
$ws.Range("M136").Value = -15450

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1751.4375
$ws.Range("I68").Value = 2129
$ws.Range("J68").Value = 1579.8182
$ws.Range("K68").Value = 6387
$ws.Range("L68").Value = 4739.4546
$ws.Range("M68").Value = -5576
$ws.Range("N68").Value = -6361.4546
$ws.Range("H71").Value = 1751.4375
$ws.Range("I71").Value = 2129
$ws.Range("J71").Value = 1579.8182
$ws.Range("K71").Value = 19161
$ws.Range("L71").Value = 14218.3638
$ws.Range("M71").Value = -15105
$ws.Range("N71").Value = -22330.3638
$ws.Range("H80").Value = 5808.3
$ws.Range("I80").Value = 5680.5
$ws.Range("K80").Value = 17041.5
$ws.Range("M80").Value = -16105.5
$ws.Range("H83").Value = 5808.3
$ws.Range("I83").Value = 5680.5
$ws.Range("K83").Value = 51124.5
$ws.Range("M83").Value = -46444.5
$ws.Range("H107").Value = 1193.9231
$ws.Range("J107").Value = 1299.375
$ws.Range("L107").Value = 3898.125
$ws.Range("N107").Value = -7738.125
$ws.Range("H137").Value = 4275
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 2188.625
$ws.Range("I140").Value = 2188.625
$ws.Range("K140").Value = 6565.875
$ws.Range("M140").Value = -1385.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6955.4287
$ws.Range("I70").Value = 6895
$ws.Range("J70").Value = 6965.5
$ws.Range("K70").Value = 6895
$ws.Range("L70").Value = 6965.5
$ws.Range("M70").Value = -6625
$ws.Range("N70").Value = -7505.5
$ws.Range("H73").Value = 6955.4287
$ws.Range("I73").Value = 6895
$ws.Range("J73").Value = 6965.5
$ws.Range("K73").Value = 6895
$ws.Range("L73").Value = 6965.5
$ws.Range("M73").Value = -5959
$ws.Range("N73").Value = -8837.5
$ws.Range("H102").Value = 2074.9443
$ws.Range("I102").Value = 1011.2222
$ws.Range("J102").Value = 3138.6667
$ws.Range("K102").Value = 1011.2222
$ws.Range("L102").Value = 3138.6667
$ws.Range("M102").Value = 610.7778
$ws.Range("N102").Value = -6382.6667
$ws.Range("H118").Value = 43000
$ws.Range("J118").Value = 43000
$ws.Range("L118").Value = 43000
$ws.Range("N118").Value = -46314

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3580.8
$ws.Range("I7").Value = 3226
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 3226
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -3114
$ws.Range("N7").Value = -5224
$ws.Range("H22").Value = 1249
$ws.Range("I22").Value = 1998
$ws.Range("J22").Value = 1076.1538
$ws.Range("K22").Value = 1998
$ws.Range("L22").Value = 1076.1538
$ws.Range("M22").Value = -1703
$ws.Range("N22").Value = -1666.1538
$ws.Range("H27").Value = 1249
$ws.Range("I27").Value = 1998
$ws.Range("J27").Value = 1076.1538
$ws.Range("K27").Value = 1998
$ws.Range("L27").Value = 1076.1538
$ws.Range("M27").Value = -1891
$ws.Range("N27").Value = -1290.1538
$ws.Range("H68").Value = 3140
$ws.Range("I68").Value = 3210
$ws.Range("K68").Value = 3210
$ws.Range("M68").Value = -2461
$ws.Range("H71").Value = 3140
$ws.Range("I71").Value = 3210
$ws.Range("K71").Value = 16050
$ws.Range("M71").Value = -12306
$ws.Range("H126").Value = 3580.8
$ws.Range("I126").Value = 3226
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9678
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7208
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 5403.1665
$ws.Range("I132").Value = 3720.2222
$ws.Range("J132").Value = 6412.933
$ws.Range("K132").Value = 11160.6666
$ws.Range("L132").Value = 19238.799
$ws.Range("M132").Value = -8630.6666
$ws.Range("N132").Value = -24298.799
$ws.Range("H136").Value = 4399.4
$ws.Range("I136").Value = 4249.25
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 12747.75
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -10197.75
$ws.Range("N136").Value = -20100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7257.9443
$ws.Range("I62").Value = 3899
$ws.Range("K62").Value = 3899
$ws.Range("M62").Value = -3275
$ws.Range("H65").Value = 7257.9443
$ws.Range("I65").Value = 3899
$ws.Range("K65").Value = 19495
$ws.Range("M65").Value = -16375
$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65766
$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67652
$ws.Range("H136").Value = 53114.7
$ws.Range("I136").Value = 2709.4666
$ws.Range("J136").Value = 204330.4
$ws.Range("K136").Value = 8128.399800000001
$ws.Range("L136").Value = 612991.2
$ws.Range("M136").Value = -5578.399800000001
$ws.Range("N136").Value = -618091.2
